$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# New row order (rows 2-12), given as arrays:
# A = Mat, B = Paterno, C = Materno, D = Nombres, E = Nombre_Largo, F = Grupo, G = Reprobadas
$data = @(
    @(20330051920389, "PACHECO",    "MAZAHUA",    "TAILY",              "ELABORA ESTRATEGIAS PARA REALIZAR LAS ACTIVIDADES DE SU ÁREA",                    "2ARHV", 2),
    @(20330051920144, "DE LA ROSA", "CASTRO",     "ALONDRA",            "ELABORA ESTRATEGIAS PARA REALIZAR LAS ACTIVIDADES DE SU ÁREA",                    "2ARHV", 2),
    @(19330051920366, "CRISTOBAL",  "BRUNO",      "DANIELA",            "EVALÚA EL DESEMPEÑO DE LA ORGANIZACIÓN UTILIZANDO HERRAMIENTAS DE CALIDAD",       "4ARHV", 2),
    @(19330051920375, "HERNANDEZ",  "ANTONIO",    "MARIA GUADALUPE",    "ASISTE EN LAS ACTIVIDADES DE CAPACITACIÓN PARA EL DESARROLLO DEL CAPITAL HUMANO", "4ARHV", 2),
    @(19330051920377, "HERNANDEZ",  "FLORES",     "PERLA",              "ASISTE EN LAS ACTIVIDADES DE CAPACITACIÓN PARA EL DESARROLLO DEL CAPITAL HUMANO", "4ARHV", 2),
    @(19330051920382, "MAZAHUA",    "IXMATLAHUA", "LUCERO",             "ASISTE EN LAS ACTIVIDADES DE CAPACITACIÓN PARA EL DESARROLLO DEL CAPITAL HUMANO", "4ARHV", 2),
    @(20330051920116, "CARRERA",    "ROMANOS",    "AMARANTA DENISSE",   "ELABORA ESTRATEGIAS PARA REALIZAR LAS ACTIVIDADES DE SU ÁREA",                    "2ARHV", 1),
    @(20330051920121, "CUATRA",     "ZOPIYACTLE", "MARIA",              "ELABORA ESTRATEGIAS PARA REALIZAR LAS ACTIVIDADES DE SU ÁREA",                    "2ARHV", 1),
    @(19330051920362, "CALIHUA",    "CALIHUA",    "JOEL",               "ASISTE EN LAS ACTIVIDADES DE CAPACITACIÓN PARA EL DESARROLLO DEL CAPITAL HUMANO", "4ARHV", 1),
    @(19330051920368, "DE LA CRUZ", "DE LA CRUZ", "OFELIA",             "ASISTE EN LAS ACTIVIDADES DE CAPACITACIÓN PARA EL DESARROLLO DEL CAPITAL HUMANO", "4ARHV", 1),
    @(19330051920378, "HERNANDEZ",  "HERNANDEZ",  "DARIANA MONSERRAT",  "ASISTE EN LAS ACTIVIDADES DE CAPACITACIÓN PARA EL DESARROLLO DEL CAPITAL HUMANO", "4ARHV", 1)
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $row++
}
